$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data row for the "Skills Imperative 2035" dataset
$ws.Range("A14").Value = "Skills Imperative 2035 employment projections by industry, sector, and qualification"
$ws.Range("B14").Value = "<a href='xxx'>NEED LINK</a>"
$ws.Range("C14").Value = "2035 (16/03/23)"
$ws.Range("D14").Value = "To be announced"

# Move the selection down, matching where the author left off editing
$ws.Range("C15").Select()
